# New Test Cases Added
# Rename Sheet1 -> "contacts" and populate it with contact test data
# (headers + two identical data rows), add mailto hyperlinks on the
# Email column, and size the columns to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "contacts"

# Header row
$headers = @("Specility", "Company", "FirstName", "MiddeName", "LastName", "Email", "Cell")
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Two identical data rows
$rowData = @("Doctor", "CLX", "Dhiraj", "V", "Redekar", "dredekar@yopmail.com", "(852) 336-6654")
for ($r = 2; $r -le 3; $r++) {
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

# Email column hyperlinks (mailto:), Excel auto-applies the built-in
# "Hyperlink" style (blue, underlined) to the affected cells.
$null = $ws.Hyperlinks.Add($ws.Range("F2"), "mailto:dredekar@yopmail.com", "", "", "dredekar@yopmail.com")
$null = $ws.Hyperlinks.Add($ws.Range("F3"), "mailto:dredekar@yopmail.com", "", "", "dredekar@yopmail.com")

# Column widths sized to the new content
$ws.Columns.Item(1).ColumnWidth = 17.666666
$ws.Columns.Item(2).ColumnWidth = 25.666666
$ws.Columns.Item(3).ColumnWidth = 22
$ws.Columns.Item(4).ColumnWidth = 22.333333
$ws.Columns.Item(5).ColumnWidth = 27.5
$ws.Columns.Item(6).ColumnWidth = 22.5
$ws.Columns.Item(7).ColumnWidth = 34.833333

# Leave the selection where the original author left it
$null = $ws.Range("B8").Select()
